# Replace NA with "", add regex test for valid lines in px file
#
# The General_MD sheet had three cells (B13:B15) holding the placeholder
# note "Befolkningsstatistikregistret indeholder ..." (stored as NA in the
# originating R data.frame). These are cleared to blank/"" cells, and the
# workbook's view state is left with General_MD active and B13:B15 selected
# (mirroring what Excel leaves behind after an interactive edit there).

$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("General_MD")

# Clear the stray "NA" note text in B13:B15 -> blank cells.
$wsGeneral.Range("B13").Value = ""
$wsGeneral.Range("B14").Value = ""
$wsGeneral.Range("B15").Value = ""

# Leave the view on General_MD with B13:B15 selected, like after editing it.
$wsGeneral.Activate()
$wsGeneral.Range("B13:B15").Select()
